$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the hyperlink/path metadata is not user-scriptable via COM; skipped.

# Row 11 existing row grows taller to fit re-wrapped text (153 -> 170)
$ws.Rows.Item(11).RowHeight = 170

# Narrow column B from 107.83 chars to 68.5 chars (stored width)
$ws.Columns.Item(2).ColumnWidth = 67.6

# Append 13 new model rows (12-24) with their equations
$ws.Range("A12").Value = 'BIOMD0000000963'
$ws.Range("B12").Value = 'odes = [
    sympy.Eq(S(t).diff(t), - beta * (S*I / (1+alpha*R)),
    sympy.Eq(I(t).diff(t), beta * (S*I / (1+alpha*R) - gamma*I),
    sympy.Eq(R(t).diff(t), gamma*I),
  ]'
$ws.Range("B12").WrapText = $true
$ws.Rows.Item(12).RowHeight = 85

$ws.Range("A13").Value = 'BIOMD0000000964'
$ws.Range("B13").Value = 'odes = [
    sympy.Eq(S(t).diff(t), b - (beta_1*S*P)/(1+alpha_1*P) - (beta_2*S*(I_A,+,I_S))/(1+alpha_2*(I_A,+,I_S)) + psi*E - mu*S),
    sympy.Eq(E(t).diff(t), (beta_1*S*P)/(1+alpha_1*P) + (beta_2*S*(I_A,+,I_S))/(1+alpha_2*(I_A,+,I_S)) - psi*E - mu*E - omega*E),
    sympy.Eq(I_A(t).diff(t), (1-delta)*omega*E - (mu+sigma)*I_A - gamma_A*I_A),
    sympy.Eq(I_S(t).diff(t), delta*omega*E - (mu+sigma)*I_S - gamma_S*I_S),
    sympy.Eq(R(t).diff(t), gamma_S*I_S + gamma_A*I_A - mu*R),
    sympy.Eq(P(t).diff(t), eta_A*I_A + eta_S*I_S - mu_p*P)
  ]'
$ws.Range("B13").WrapText = $true
$ws.Rows.Item(13).RowHeight = 170

$ws.Range("A14").Value = 'BIOMD0000000970'
$ws.Range("B14").Value = 'odes = [
    sympy.Eq(S(t).diff(t), -r_1*beta_1*I*S/N - r_2*beta_2*E*S/N),
    sympy.Eq(E(t).diff(t), r_1*beta*I*S/N - alpha*E + r_2*beta_2*E*S/N),
    sympy.Eq(I(t).diff(t), alpha*E - gamma*I),
    sympy.Eq(R(t).diff(t), gamma*I)
]'
$ws.Range("B14").WrapText = $true
$ws.Rows.Item(14).RowHeight = 102

$ws.Range("A15").Value = 'BIOMD0000000971'
$ws.Range("B15").Value = 'odes = [
    sympy.Eq(S(t).diff(t), -(beta*c + c*q*(1 - beta))*S*(I + theta*A) + lambda_*S_q),
    sympy.Eq(E(t).diff(t), beta*c*(1 - q)*S*(I + theta*A) - sigma*E),
    sympy.Eq(I(t).diff(t), sigma*rho*E - (delta_I + alpha + gamma_I)*I),
    sympy.Eq(A(t).diff(t), sigma*(1 - rho)*E - gamma_A*A),
    sympy.Eq(S_q(t).diff(t), (1 - beta)*c*q*S*(I + theta*A) - lambda_*S_q),
    sympy.Eq(E_q(t).diff(t), beta*c*q*S*(I + theta*A) - delta_q*E_q),
    sympy.Eq(H(t).diff(t), delta_I*I + delta_q*E_q - (alpha + gamma_H)*H),
    sympy.Eq(R(t).diff(t), gamma_I*I + gamma_A*A + gamma_H*H)
]'
$ws.Range("B15").WrapText = $true
$ws.Rows.Item(15).RowHeight = 170

$ws.Range("A16").Value = 'BIOMD0000000972'
$ws.Range("B16").Value = 'odes = [
    sympy.Eq(S(t).diff(t), -(beta*c(t) + c(t)*q*(1 - beta))*S*(I + theta*A) + lambda_*S_q),
    sympy.Eq(E(t).diff(t), beta*c(t)*(1 - q)*S*(I + theta*A) - sigma*E),
    sympy.Eq(I(t).diff(t), sigma*rho*E - (delta_I(t) + alpha + gamma_I)*I),
    sympy.Eq(A(t).diff(t), sigma*(1 - rho)*E - gamma_A*A),
    sympy.Eq(S_q(t).diff(t), (1 - beta)*c(t)*q*S*(I + theta*A) - lambda_*S_q),
    sympy.Eq(E_q(t).diff(t), beta*c(t)*q*S*(I + theta*A) - delta_q*E_q),
    sympy.Eq(H(t).diff(t), delta_I(t)*I + delta_q*E_q - (alpha + gamma_H)*H),
    sympy.Eq(R(t).diff(t), gamma_I*I + gamma_A*A + gamma_H*H)
]'
$ws.Range("B16").WrapText = $true
$ws.Rows.Item(16).RowHeight = 187

$ws.Range("A17").Value = 'BIOMD0000000974'
$ws.Range("B17").Value = 'odes = [
    sympy.Eq(S(t).diff(t), Lambda - mu*S - beta*S*I/N),
    sympy.Eq(E(t).diff(t), beta*S*I/N - (mu + epsilon)*E),
    sympy.Eq(I(t).diff(t), epsilon*E - (gamma + mu + alpha)*I),
    sympy.Eq(R(t).diff(t), gamma*I - mu*R)
]'
$ws.Range("B17").WrapText = $true
$ws.Rows.Item(17).RowHeight = 102

$ws.Range("A18").Value = 'BIOMD0000000976'
$ws.Range("B18").Value = 'odes = [
    sympy.Eq(S(t).diff(t), -tau(t)*S(t)*(I_1(t) + I_2(t))),
    sympy.Eq(I_1(t).diff(t), alpha*tau(t)*S(t)*(I_1(t) + I_2(t)) - gamma_1*I_1(t)),
    sympy.Eq(I_2(t).diff(t), (1 - alpha)*tau(t)*S(t)*(I_1(t) + I_2(t)) - (gamma_2 + mu)*I_2(t)),
    sympy.Eq(R(t).diff(t), gamma_1*I_1(t) + gamma_2*I_2(t))
]'
$ws.Range("B18").WrapText = $true
$ws.Rows.Item(18).RowHeight = 119

$ws.Range("A19").Value = 'BIOMD0000000977'
$ws.Range("B19").Value = 'odes = [
    sympy.Eq(S(t).diff(t), Lambda_s - (beta_s + rho_s*(1 - beta_s))*epsilon_s*S*I/N - delta*S + m_s*S_q),
    sympy.Eq(S_q(t).diff(t), (1 - beta_s)*epsilon_s*rho_s*S*I/N - (m_s + delta)*S_q),
    sympy.Eq(A(t).diff(t), beta_s*(1 - rho_s)*epsilon_s*S*I/N - (gamma_a + xi_a + delta)*A),
    sympy.Eq(I(t).diff(t), gamma_a*A - (gamma_i + xi_i + delta)*I),
    sympy.Eq(I_q(t).diff(t), beta_s*epsilon_s*rho_s*S*I/N + gamma_i*I - (xi_q + delta)*I_q),
    sympy.Eq(R(t).diff(t), xi_a*A + xi_i*I + xi_q*I_q - delta*R)
]'
$ws.Range("B19").WrapText = $true
$ws.Rows.Item(19).RowHeight = 187

$ws.Range("A20").Value = 'BIOMD0000000978'
$ws.Range("B20").Value = 'odes = [
    sympy.Eq(S(t).diff(t), -(1 - epsilon)*beta*S*I/N),
    sympy.Eq(E(t).diff(t), (1 - epsilon)*beta*S*I/N - sigma*E),
    sympy.Eq(I(t).diff(t), sigma*E - gamma*I),
    sympy.Eq(R(t).diff(t), gamma*I)
]'
$ws.Range("B20").WrapText = $true
$ws.Rows.Item(20).RowHeight = 102

$ws.Range("A21").Value = 'BIOMD0000000979'
$ws.Range("B21").Value = 'odes = [
    sympy.Eq(S(t).diff(t), -beta(t)*S(t)/N*I(t) + omega*R(t)),
    sympy.Eq(E(t).diff(t), beta(t)*S(t)*/N*I(t) - sigma*E(t)),
    sympy.Eq(I(t).diff(t), sigma*E(t) - gamma*I(t)),
    sympy.Eq(R(t).diff(t), gamma*I(t) - omega*R(t))
]'
$ws.Range("B21").WrapText = $true
$ws.Rows.Item(21).RowHeight = 102

$ws.Range("A22").Value = 'BIOMD0000000983'
$ws.Range("B22").Value = 'odes = [
    sympy.Eq(S_c(t).diff(t), m(t)*S_u - (1 - m(t))*S_c),
    sympy.Eq(S_u(t).diff(t), (1 - m(t))*S_c - m(t)*S_u - beta*(n*I_r + I_u)*S_u + theta*(1 - lambda_)*Q),
    sympy.Eq(E(t).diff(t), (1 - sigma)*beta*(n*I_r + I_u)*S_u - mu*E),
    sympy.Eq(I_r(t).diff(t), mu*f*E + theta*lambda*Q - eta_r*I_r),
    sympy.Eq(I_u(t).diff(t), mu*(1 - f)*E - eta_u*I_u),
    sympy.Eq(R(t).diff(t), eta_r*q*I_r + eta_u*I_u),
    sympy.Eq(Q(t).diff(t), sigma*beta*(n*I_r + I_u)*S_u - theta*Q)
]'
$ws.Range("B22").WrapText = $true
$ws.Rows.Item(22).RowHeight = 170

$ws.Range("A23").Value = 'BIOMD0000000984'
$ws.Range("B23").Value = 'odes = [
    sympy.Eq(S(t).diff(t), -beta*S/N*I),
    sympy.Eq(E(t).diff(t), beta*S/N*I - omega*E),
    sympy.Eq(I(t).diff(t), omega*E - gamma*I),
    sympy.Eq(R(t).diff(t), gamma*I)
]'
$ws.Range("B23").WrapText = $true
$ws.Rows.Item(23).RowHeight = 102

$ws.Range("A24").Value = 'BIOMD0000000991A'
$ws.Range("B24").Value = 'odes = [
    sympy.Eq(S(t).diff(t), -beta_c*(alpha*A + I)*S/(N_h - I_D)),
    sympy.Eq(E(t).diff(t), beta_c*(alpha*A + I)*S/(N_h - I_D) - sigma*E),
    sympy.Eq(A(t).diff(t), nu*sigma*E - (theta + gamma_a)*A),
    sympy.Eq(I(t).diff(t), (1 - nu)*sigma*E - (psi + gamma_O + d_O)*I),
    sympy.Eq(I_D(t).diff(t), theta*A + psi*I - (gamma_i + d_D)*I_D),
    sympy.Eq(R(t).diff(t), gamma_i*I_D + gamma_a*A + gamma_O*I)
]'
$ws.Range("B24").WrapText = $true
$ws.Rows.Item(24).RowHeight = 136

# Match the final view state: zoomed in, scrolled near the bottom, selection on the next empty row
$excel.ActiveWindow.Zoom = 200
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B25").Select()
